$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42 used to be a (mostly) blank template row with a leftover border
# style on F42 - clear it out first so the new data doesn't inherit that
# stray formatting.
$ws.Range("A42:N42").Clear()

# Add the new "GrenadierWeapon" row (row 42)
$ws.Range("A42").Value = 9004
$ws.Range("B42").Value = "GrenadierWeapon"
$ws.Range("C42").Value = "Common"
$ws.Range("D42").Value = 15
$ws.Range("E42").Value = "몬스터가 쓰는 무기인 것 같다 버리고 가자."
$ws.Range("F42").Value = 0.5
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 7
$ws.Range("I42").Value = $true
$ws.Range("J42").Value = 0.7
$ws.Range("K42").Value = 0.5
$ws.Range("L42").Value = 97
$ws.Range("M42").Value = $false
$ws.Range("N42").Value = 0

# Move the active selection to B42
$ws.Range("B42").Select() | Out-Null
